# static_params.xlsx - minor fixes and features
#
# - EQUITY row: subaccount placeholder "SysPerp" -> "debug"
# - SLIPPAGE_OVERRIDE: 0.0002 -> 0 (disabled)
# - CONCENTRATION_LIMIT: 0.75 -> 2, with a new comment explaining the unit
# - cosmetic: last-used cell + column widths refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# EQUITY (row 5): value "SysPerp" -> "debug"
$ws.Range("B5").Value = "debug"

# SLIPPAGE_OVERRIDE (row 9): 0.0002 -> 0
$ws.Range("B9").Value = 0

# CONCENTRATION_LIMIT (row 10): 0.75 -> 2, and document the unit in C10
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "max volume/open interest share"

# Cosmetic: refresh column widths (A/B tightened) and restore last selection
$ws.Columns.Item(1).ColumnWidth = 25.85
$ws.Columns.Item(2).ColumnWidth = 10.17

$ws.Range("D9").Select()
